# Inicialización de la base de datos
# Adds employee record #8 (Luis Jose Campos Perez) as a new row (row 9)
# at the bottom of the "Reporte_Empleados" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Luis"
$ws.Range("C9").Value = "Jose"
$ws.Range("D9").Value = "Campos"
$ws.Range("E9").Value = "Perez"
# Leading apostrophe forces text storage so numeric-looking values
# ("17342541", "0101.2222222") and the date-looking value
# ("11-10-2013") are kept as literal text, matching the source data.
$ws.Range("F9").Value = "'17342541"
$ws.Range("G9").Value = "mari@hma,com"
$ws.Range("H9").Value = "mama"
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = "Masculino"
$ws.Range("L9").Value = "Soltero"
$ws.Range("M9").Value = "COntador"
$ws.Range("N9").Value = "Finanzas"
$ws.Range("O9").Value = "'11-10-2013"
$ws.Range("P9").Value = "Caracas"
$ws.Range("Q9").Value = "Quincenal"
$ws.Range("R9").Value = "Activo"
$ws.Range("S9").Value = "BBVA"
$ws.Range("T9").Value = "'0101.2222222"
